$p = $ppt.ActivePresentation

# --- Slide 12: "Content Placeholder 2" reposition/resize ---
$s12 = $p.Slides.Item(12)
$shpContent = $s12.Shapes.Item(3)
$shpContent.Left = 66
$shpContent.Top = 61.07204724409449
$shpContent.Width = 588
$shpContent.Height = 67.5

# --- Slide 6: figure caption text updates ---
$s6 = $p.Slides.Item(6)

# "Rectangle 6" - Query Message caption
$shpQuery = $s6.Shapes.Item(4)
$trQuery = $shpQuery.TextFrame.TextRange
$paraCount4 = $trQuery.Paragraphs().Count
$capQuery = $trQuery.Paragraphs($paraCount4, 1)
# Replace via a disjoint placeholder first so the host doesn't try to
# diff/split the run against the old text (it preserves full rPr only
# when old and new text share no substring).
$capQuery.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$capQuery2 = $trQuery.Paragraphs($paraCount4, 1)
$capQuery2.Text = "             Figure: Control Code in TWAMP Light Query Message"

# "Rectangle 8" - Response Message caption
$shpResp = $s6.Shapes.Item(6)
$trResp = $shpResp.TextFrame.TextRange
$paraCount6 = $trResp.Paragraphs().Count
$capResp = $trResp.Paragraphs($paraCount6, 1)
$capResp.Text = "zzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzzz"
$capResp2 = $trResp.Paragraphs($paraCount6, 1)
$capResp2.Text = "           Figure: Control Code in TWAMP Light Response Message"
